$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("steel")

$ws.Range("A34").Value = "china-HeEtAl2017"
$ws.Range("B34").Value = "coal"
$ws.Range("C34").Value = "air"
$ws.Range("D34").Value = "BF+BOF"
$ws.Range("D34").WrapText = $true
